# Update "想去人数" (interested-count) figures on both the "展览" sheet
# and the "全部类型" sheet to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 12163
$ws1.Range("F5").Value = 4481
$ws1.Range("F10").Value = 2595
$ws1.Range("F11").Value = 1125
$ws1.Range("F12").Value = 190
$ws1.Range("F13").Value = 64
$ws1.Range("F14").Value = 5268
$ws1.Range("F17").Value = 548
$ws1.Range("F18").Value = 11443
$ws1.Range("F19").Value = 11510

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12163
$ws4.Range("F5").Value = 4481
$ws4.Range("F10").Value = 2595
$ws4.Range("F12").Value = 1125
$ws4.Range("F13").Value = 190
$ws4.Range("F14").Value = 64
$ws4.Range("F15").Value = 5268
$ws4.Range("F18").Value = 548
$ws4.Range("F19").Value = 11443
$ws4.Range("F20").Value = 11510
